# Stand Up Meeting Week 9 - fill in Guillermo's row (16-18) with his
# own status-report text, and widen column C so the longer text fits.
#
# NOTE on ordering: new shared strings are appended to the shared-string
# table in the order cells are first assigned a *new* (not-yet-seen)
# string value. To reproduce the exact target shared-string indices we
# must assign C17 before D16 before D17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New text (first occurrences, order matters) ---------------------
$ws.Range("C17").Value = "Reunión para asignar tareas, revisión a lo que hice"
$ws.Range("D16").Value = "Asistimos a la reunión y se asignaron tareas en la tabla de casos de uso. Andrea y yo terminamos el modelo de clases, le fata la cardinalidad"
$ws.Range("D17").Value = "nada"

# --- Remaining cells, reusing existing shared strings -----------------
$ws.Range("C16").Value = "Nada"
$ws.Range("E16").Value = "Trabajamos en el diagrama de clases."
$ws.Range("F16").Value = "Nada"
$ws.Range("G16").Value = "Tener en cuenta las correcciones"

$ws.Range("E17").Value = "Nada"
$ws.Range("F17").Value = "Mostrar los avances en la clase, anotar las correcciones"
$ws.Range("G17").Value = "Asistir a la reunión"

$ws.Range("C18").Value = "Nada"
$ws.Range("D18").Value = "Nada"
$ws.Range("E18").Value = "Nada"
$ws.Range("F18").Value = "Nada"
$ws.Range("G18").Value = "Nada"

# --- Widen column C to fit the longer text -----------------------------
$ws.Columns.Item(3).ColumnWidth = 58.88671875

# --- Scroll / selection state, matching the saved view in the workbook -
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G18").Select() | Out-Null
